$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.994735717773438
$ws.Range("B1").Value = 5.216355800628662
$ws.Range("C1").Value = 2.918453931808472
$ws.Range("D1").Value = 1.179003357887268
$ws.Range("E1").Value = 0.8124327659606934
